$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price (D) and volume-change (E) values.
# D-column cells are forced to Text format before assignment so that
# Excel does not reinterpret numeric-looking strings (e.g. "1.002") as
# numbers; the style is then reset to Normal so no stray cell formatting
# is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.844.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.707.04'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.48%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3943'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4049'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.521'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.002'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.47'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08927'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.305'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.74'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.041'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001329'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.727.78'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '100.61'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07043'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.79'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.094'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("E23").Value = '  +2.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.822.19'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.82%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.240'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.363'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.88'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.82%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.57'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.464'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +10.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '136.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.173'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08947'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.594'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.94%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.089'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.18'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.982'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2764'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '14.55'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09223'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02773'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.463'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7735'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '15.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.93%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7226'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.581'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.33%  '
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.002'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '140.72'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.91%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.326'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '91.42'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.08013'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.38%  '
